# Update predicted-signal and actual-return columns (AC, AD) on rows 3-6
# to reflect the recalculated values produced after the test/train split
# was refactored into smaller modules.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC3").Value = 1

$ws.Range("AC4").Value = 1
$ws.Range("AD4").Value = -0.002178269582643555

$ws.Range("AC5").Value = 1
$ws.Range("AD5").Value = -0.006505413901501833

$ws.Range("AD6").Value = 0.01647989452867504
